# ProductHistorical_TestData.xlsx - "makes relative xpath in productHistorical package"
#
# The test-data generator's "Web Data" / "Child Web Data" run id moved from
# 64 to 33, and the CreateOrder sheet's product names / order totals were
# reworked (product names now have a space before the number, and the order
# totals are now stored as the already-formatted display string instead of
# the bare numeric amount).

$wb = $excel.ActiveWorkbook

$nbsp = [char]0x00A0
$euro = [char]0x20AC

# ---------------------------------------------------------------------
# 1) LoginSignup
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LoginSignup")
$ws.Range("C1").Value = "Web Data 33"
$ws.Range("D1").Value = "Web Data 33"
$ws.Range("N1").Value = "Child Web Data 33"
$ws.Range("U1").Value = "Successfully created Child Web Data 33. You can now login with the username admin after your password is set. Password reset link is sent to your email."
$ws.Range("D1").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) AddCurrency
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AddCurrency")
$ws.Range("C1").Value = "Web Data 33"
$ws.Range("D1").Value = "Child Web Data 33"
$ws.Range("F1").Value = "Working as admin Child Web Data 33 X"
$ws.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) AddProductCategory1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AddProductCategory1")
$ws.Range("C1").Value = "Web Data 33"
$ws.Range("E1").Value = "Child Web Data 33"
$ws.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------------
# 4) AddCustomer
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AddCustomer")
$ws.Range("C1").Value = "Web Data 33"
$ws.Range("C2").Value = "Child Web Data 33"
$ws.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------------
# 5) CreateOrder
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CreateOrder")

$ws.Range("C1").Value = "Web Data 33"

$ws.Range("F1").Value = "Test Product 1"
$ws.Range("G1").Value = "Test Product 1"
$ws.Range("K1").Value = "Total = C`$123.00"
$ws.Range("L1").Value = "Total = " + $euro + "96.00"
$ws.Range("M1").Value = "Total = " + $euro + "126.00"

$ws.Range("C2").Value = "Child Web Data 33"
$ws.Range("F2").Value = "Test Product 2"
$ws.Range("G2").Value = "Test Product 3"
$ws.Range("K2").Value = "Total = C`$123.00"
$ws.Range("L2").Value = "Total = " + $euro + "96.00"
$ws.Range("M2").Value = "Total = " + $euro + "126.00"

$ws.Range("C3").Value = "Working as admin Child Web Data 33" + $nbsp + "X"
$ws.Range("F3").Value = "Test Product 4"
$ws.Range("G3").Value = "Test Product 4"
$ws.Range("K3").Value = "Total = C`$132.00"
$ws.Range("L3").Value = "Total = " + $euro + "105.00"
$ws.Range("M3").Value = "Total = " + $euro + "135.00"

$ws.Range("K4").Value = "Total = C`$141.00"
$ws.Range("L4").Value = "Total = " + $euro + "114.00"
$ws.Range("M4").Value = "Total = " + $euro + "144.00"

$ws.Range("M5").Select() | Out-Null
